$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("class_schedule")

# NOTE: new shared strings are appended to the sst in first-write order, so
# write D11's new text before D6's to reproduce the author's original
# sharedStrings.xml ordering (JVP/Views-and-Copies string first, then the
# Follow-this-link/Python-v-R string).

# Both replacement texts start with a literal "- " (hyphen, like the text
# they replace), so set them through .Formula with a leading apostrophe -
# Excel's plain-text escape - to keep the cells' existing quote-prefixed
# style (xf s="4") instead of having a plain .Value assignment silently
# reformat them to the un-prefixed style used elsewhere in the sheet.

# D11 (topic: "Pandas: Indices & Missing"): previously just "- JVP pp 115-139";
# now also links the new "Views and Copies in Pandas" notebook.
$newD11 = "'- JVP pp 115-139`n- ``Views and Copies in Pandas <views_and_copies_in_pandas.ipynb>```_"
$ws.Range("D11").Formula = $newD11

# D6 (topic: "Python v. R / variables as pointers"): previously just the
# "Python v R: Pointers" notebook link; now also points students to the
# Ipython-in-depth binder notebook first.
$newD6 = "'- ``Follow this link <https://gke.mybinder.org/v2/gh/ipython/ipython-in-depth/master?filepath=binder/Index.ipynb>```_ , then click ""Ipython - Beyond plain python"" and read that notebook. `n- ``Python v R: Pointers <python_v_r.ipynb>```_"
$ws.Range("D6").Formula = $newD6

# Row heights grew to fit the new, longer wrapped text.
$ws.Rows.Item(6).RowHeight = 71
$ws.Rows.Item(11).RowHeight = 43

# Move the active selection from D25 to D6.
$ws.Activate()
$ws.Range("D6").Select() | Out-Null
